$wb = $excel.ActiveWorkbook

# Sheets: 1=NewLoanInput, 2=Summary, 3=Repayment schedule, 4=Transactions
$wsRepay = $wb.Worksheets.Item(3)
$wsTrans = $wb.Worksheets.Item(4)

# Insert a new (blank) column before column N ("Late") on the
# "Repayment schedule" sheet - this shifts old N/O/P -> O/P/Q.
$wsRepay.Columns("N:N").Insert()

# Match the width of the newly inserted column to the target layout
# (stored column width of 11, not a bestFit column like its neighbours).
$wsRepay.Columns("N:N").ColumnWidth = 10.2

# Update the selected / active sheet + cell.
# "Transactions" loses the tab-selected / active-cell state ...
$wsTrans.Activate()
$wsTrans.Range("D22").Select()

# ... which moves to "Repayment schedule".
$wsRepay.Activate()
$wsRepay.Range("S6").Select()
